# Natmi following Dr Hou advice
# Recomputed NATMI ligand/receptor edge-weight statistics for Jag2-Notch2
# after bumping "expressing cells" counts (E/K columns) from 1 to 3 for
# every data row (rows 2-17), which cascades into new average/total
# expression, specificity and edge-weight values (columns G-J, M-T).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 20.70050833333333
$ws.Cells.Item(2, 8).Value = 62.101525
$ws.Cells.Item(2, 9).Value = 0.8277101186170105
$ws.Cells.Item(2, 10).Value = 0.8277101186170105
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 2.993142333333334
$ws.Cells.Item(2, 14).Value = 8.979427000000001
$ws.Cells.Item(2, 15).Value = 0.03484385887642424
$ws.Cells.Item(2, 16).Value = 0.03484385887642424
$ws.Cells.Item(2, 17).Value = 61.95956781401944
$ws.Cells.Item(2, 18).Value = 557.6361103261751
$ws.Cells.Item(2, 19).Value = 0.02884061456367948
$ws.Cells.Item(2, 20).Value = 0.02884061456367949
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 20.70050833333333
$ws.Cells.Item(3, 8).Value = 62.101525
$ws.Cells.Item(3, 9).Value = 0.8277101186170105
$ws.Cells.Item(3, 10).Value = 0.8277101186170105
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 31.995262
$ws.Cells.Item(3, 14).Value = 95.985786
$ws.Cells.Item(3, 15).Value = 0.3724642097459734
$ws.Cells.Item(3, 16).Value = 0.3724642097459735
$ws.Cells.Item(3, 17).Value = 662.3181876581833
$ws.Cells.Item(3, 18).Value = 5960.86368892365
$ws.Cells.Item(3, 19).Value = 0.3082923952294307
$ws.Cells.Item(3, 20).Value = 0.3082923952294308
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 20.70050833333333
$ws.Cells.Item(4, 8).Value = 62.101525
$ws.Cells.Item(4, 9).Value = 0.8277101186170105
$ws.Cells.Item(4, 10).Value = 0.8277101186170105
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 34.28929533333334
$ws.Cells.Item(4, 14).Value = 102.867886
$ws.Cells.Item(4, 15).Value = 0.3991695798295478
$ws.Cells.Item(4, 16).Value = 0.3991695798295478
$ws.Cells.Item(4, 17).Value = 709.8058437917945
$ws.Cells.Item(4, 18).Value = 6388.252594126149
$ws.Cells.Item(4, 19).Value = 0.3303967002690172
$ws.Cells.Item(4, 20).Value = 0.3303967002690172
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 20.70050833333333
$ws.Cells.Item(5, 8).Value = 62.101525
$ws.Cells.Item(5, 9).Value = 0.8277101186170105
$ws.Cells.Item(5, 10).Value = 0.8277101186170105
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 16.62387466666667
$ws.Cells.Item(5, 14).Value = 49.871624
$ws.Cells.Item(5, 15).Value = 0.1935223515480544
$ws.Cells.Item(5, 16).Value = 0.1935223515480545
$ws.Cells.Item(5, 17).Value = 344.1226560696222
$ws.Cells.Item(5, 18).Value = 3097.1039046266
$ws.Cells.Item(5, 19).Value = 0.1601804085548829
$ws.Cells.Item(5, 20).Value = 0.160180408554883
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 0.7925996666666667
$ws.Cells.Item(6, 8).Value = 2.377799
$ws.Cells.Item(6, 9).Value = 0.03169210888681734
$ws.Cells.Item(6, 10).Value = 0.03169210888681734
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 2.993142333333334
$ws.Cells.Item(6, 14).Value = 8.979427000000001
$ws.Cells.Item(6, 15).Value = 0.03484385887642424
$ws.Cells.Item(6, 16).Value = 0.03484385887642424
$ws.Cells.Item(6, 17).Value = 2.372363615685889
$ws.Cells.Item(6, 18).Value = 21.351272541173
$ws.Cells.Item(6, 19).Value = 0.001104275369548534
$ws.Cells.Item(6, 20).Value = 0.001104275369548534
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 0.7925996666666667
$ws.Cells.Item(7, 8).Value = 2.377799
$ws.Cells.Item(7, 9).Value = 0.03169210888681734
$ws.Cells.Item(7, 10).Value = 0.03169210888681734
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 31.995262
$ws.Cells.Item(7, 14).Value = 95.985786
$ws.Cells.Item(7, 15).Value = 0.3724642097459734
$ws.Cells.Item(7, 16).Value = 0.3724642097459735
$ws.Cells.Item(7, 17).Value = 25.35943399611267
$ws.Cells.Item(7, 18).Value = 228.234905965014
$ws.Cells.Item(7, 19).Value = 0.01180417629171176
$ws.Cells.Item(7, 20).Value = 0.01180417629171176
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 0.7925996666666667
$ws.Cells.Item(8, 8).Value = 2.377799
$ws.Cells.Item(8, 9).Value = 0.03169210888681734
$ws.Cells.Item(8, 10).Value = 0.03169210888681734
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 34.28929533333334
$ws.Cells.Item(8, 14).Value = 102.867886
$ws.Cells.Item(8, 15).Value = 0.3991695798295478
$ws.Cells.Item(8, 16).Value = 0.3991695798295478
$ws.Cells.Item(8, 17).Value = 27.17768405143489
$ws.Cells.Item(8, 18).Value = 244.599156462914
$ws.Cells.Item(8, 19).Value = 0.01265052578826315
$ws.Cells.Item(8, 20).Value = 0.01265052578826315
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 0.7925996666666667
$ws.Cells.Item(9, 8).Value = 2.377799
$ws.Cells.Item(9, 9).Value = 0.03169210888681734
$ws.Cells.Item(9, 10).Value = 0.03169210888681734
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 16.62387466666667
$ws.Cells.Item(9, 14).Value = 49.871624
$ws.Cells.Item(9, 15).Value = 0.1935223515480544
$ws.Cells.Item(9, 16).Value = 0.1935223515480545
$ws.Cells.Item(9, 17).Value = 13.17607751950844
$ws.Cells.Item(9, 18).Value = 118.584697675576
$ws.Cells.Item(9, 19).Value = 0.006133131437293885
$ws.Cells.Item(9, 20).Value = 0.006133131437293886
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 1.562510666666667
$ws.Cells.Item(10, 8).Value = 4.687532
$ws.Cells.Item(10, 9).Value = 0.06247701111592723
$ws.Cells.Item(10, 10).Value = 0.06247701111592723
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 2.993142333333334
$ws.Cells.Item(10, 14).Value = 8.979427000000001
$ws.Cells.Item(10, 15).Value = 0.03484385887642424
$ws.Cells.Item(10, 16).Value = 0.03484385887642424
$ws.Cells.Item(10, 17).Value = 4.676816822684889
$ws.Cells.Item(10, 18).Value = 42.09135140416401
$ws.Cells.Item(10, 19).Value = 0.002176940158344156
$ws.Cells.Item(10, 20).Value = 0.002176940158344157
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 1.562510666666667
$ws.Cells.Item(11, 8).Value = 4.687532
$ws.Cells.Item(11, 9).Value = 0.06247701111592723
$ws.Cells.Item(11, 10).Value = 0.06247701111592723
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 31.995262
$ws.Cells.Item(11, 14).Value = 95.985786
$ws.Cells.Item(11, 15).Value = 0.3724642097459734
$ws.Cells.Item(11, 16).Value = 0.3724642097459735
$ws.Cells.Item(11, 17).Value = 49.99293815779467
$ws.Cells.Item(11, 18).Value = 449.936443420152
$ws.Cells.Item(11, 19).Value = 0.02327045057258423
$ws.Cells.Item(11, 20).Value = 0.02327045057258424
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 1.562510666666667
$ws.Cells.Item(12, 8).Value = 4.687532
$ws.Cells.Item(12, 9).Value = 0.06247701111592723
$ws.Cells.Item(12, 10).Value = 0.06247701111592723
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 34.28929533333334
$ws.Cells.Item(12, 14).Value = 102.867886
$ws.Cells.Item(12, 15).Value = 0.3991695798295478
$ws.Cells.Item(12, 16).Value = 0.3991695798295478
$ws.Cells.Item(12, 17).Value = 53.57738971081689
$ws.Cells.Item(12, 18).Value = 482.196507397352
$ws.Cells.Item(12, 19).Value = 0.02493892227615066
$ws.Cells.Item(12, 20).Value = 0.02493892227615066
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 1.562510666666667
$ws.Cells.Item(13, 8).Value = 4.687532
$ws.Cells.Item(13, 9).Value = 0.06247701111592723
$ws.Cells.Item(13, 10).Value = 0.06247701111592723
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 16.62387466666667
$ws.Cells.Item(13, 14).Value = 49.871624
$ws.Cells.Item(13, 15).Value = 0.1935223515480544
$ws.Cells.Item(13, 16).Value = 0.1935223515480545
$ws.Cells.Item(13, 17).Value = 25.97498148799644
$ws.Cells.Item(13, 18).Value = 233.774833391968
$ws.Cells.Item(13, 19).Value = 0.01209069810884817
$ws.Cells.Item(13, 20).Value = 0.01209069810884818
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 1.953751
$ws.Cells.Item(14, 8).Value = 5.861253
$ws.Cells.Item(14, 9).Value = 0.07812076138024482
$ws.Cells.Item(14, 10).Value = 0.07812076138024483
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 2.993142333333334
$ws.Cells.Item(14, 14).Value = 8.979427000000001
$ws.Cells.Item(14, 15).Value = 0.03484385887642424
$ws.Cells.Item(14, 16).Value = 0.03484385887642424
$ws.Cells.Item(14, 17).Value = 5.847854826892333
$ws.Cells.Item(14, 18).Value = 52.630693442031
$ws.Cells.Item(14, 19).Value = 0.002722028784852063
$ws.Cells.Item(14, 20).Value = 0.002722028784852064
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 1.953751
$ws.Cells.Item(15, 8).Value = 5.861253
$ws.Cells.Item(15, 9).Value = 0.07812076138024482
$ws.Cells.Item(15, 10).Value = 0.07812076138024483
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 31.995262
$ws.Cells.Item(15, 14).Value = 95.985786
$ws.Cells.Item(15, 15).Value = 0.3724642097459734
$ws.Cells.Item(15, 16).Value = 0.3724642097459735
$ws.Cells.Item(15, 17).Value = 62.510775127762
$ws.Cells.Item(15, 18).Value = 562.5969761498579
$ws.Cells.Item(15, 19).Value = 0.02909718765224665
$ws.Cells.Item(15, 20).Value = 0.02909718765224666
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 1.953751
$ws.Cells.Item(16, 8).Value = 5.861253
$ws.Cells.Item(16, 9).Value = 0.07812076138024482
$ws.Cells.Item(16, 10).Value = 0.07812076138024483
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 34.28929533333334
$ws.Cells.Item(16, 14).Value = 102.867886
$ws.Cells.Item(16, 15).Value = 0.3991695798295478
$ws.Cells.Item(16, 16).Value = 0.3991695798295478
$ws.Cells.Item(16, 17).Value = 66.99274504679533
$ws.Cells.Item(16, 18).Value = 602.9347054211579
$ws.Cells.Item(16, 19).Value = 0.03118343149611669
$ws.Cells.Item(16, 20).Value = 0.03118343149611669
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 1.953751
$ws.Cells.Item(17, 8).Value = 5.861253
$ws.Cells.Item(17, 9).Value = 0.07812076138024482
$ws.Cells.Item(17, 10).Value = 0.07812076138024483
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 16.62387466666667
$ws.Cells.Item(17, 14).Value = 49.871624
$ws.Cells.Item(17, 15).Value = 0.1935223515480544
$ws.Cells.Item(17, 16).Value = 0.1935223515480545
$ws.Cells.Item(17, 17).Value = 25.97498148799644
$ws.Cells.Item(17, 18).Value = 233.774833391968
$ws.Cells.Item(17, 19).Value = 0.01209069810884817
$ws.Cells.Item(17, 20).Value = 0.01209069810884818
